# Insert a new data row right before the current row 248, shifting all
# subsequent rows (old 248-308) down by one (they become new rows 249-309).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(248).Insert()

# Populate the newly inserted row 248 with the new market record.
$ws.Cells.Item(248, 1).Value = 5
$ws.Cells.Item(248, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(248, 3).Value = "Maule"
$ws.Cells.Item(248, 4).Value = Get-Date -Year 2022 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(248, 5).Value = 7
$ws.Cells.Item(248, 6).Value = "Fruta"
$ws.Cells.Item(248, 7).Value = 100102
$ws.Cells.Item(248, 8).Value = "Cítricos"
$ws.Cells.Item(248, 9).Value = 100102004
$ws.Cells.Item(248, 10).Value = "Mandarina"
$ws.Cells.Item(248, 11).Value = "Murcott"
$ws.Cells.Item(248, 12).Value = "Primera"
$ws.Cells.Item(248, 13).Value = 250
$ws.Cells.Item(248, 14).Value = 7000
$ws.Cells.Item(248, 15).Value = 7000
$ws.Cells.Item(248, 16).Value = 7000
$ws.Cells.Item(248, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(248, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(248, 19).Value = 389
$ws.Cells.Item(248, 20).Value = 18
